$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typos / renumbered text in row 3 (existing risk: Laptop Failure) ---
$ws.Range("F3").Value = "Continue development on my own laptop until a replacement arrives "
$ws.Range("G3").Value = "Make sure laptop is update and working correctly regulary "
$ws.Range("H3").Value = "No updates at of 23/04/2021. "

# --- Row 4: Passwords / credentials pushed to github ---
$ws.Range("A4").Value = "Critial information pushed to github."
$ws.Range("B4").Value = "Results in exposure of passwords " + [char]10 + "and other information to individuals who should not be able to see it"
$ws.Range("C4").Value = "low"
$ws.Range("D4").Value = "high"
$ws.Range("E4").Value = "Sean Kaila"
$ws.Range("F4").Value = "if the information where to reach the wrong had," + [char]10 + " change the passwords and other informationassociated with data breach"
$ws.Range("G4").Value = "Use a GitIgnore file to hide connection strings, " + [char]10 + "passwords and test result files from being pushed to github"
$ws.Range("H4").Value = "No updates as of 08/05/2021"
$ws.Range("B4").WrapText = $true
$ws.Range("F4").WrapText = $true
$ws.Range("G4").WrapText = $true

# --- Row 5: code Loss (Local) ---
$ws.Range("A5").Value = "code Loss (Local)"
$ws.Range("B5").Value = "Key progress of project lost. Therefore " + [char]10 + "not able to deliver project deliverables"
$ws.Range("C5").Value = "low"
$ws.Range("D5").Value = "low"
$ws.Range("E5").Value = "Sean Kaila"
$ws.Range("F5").Value = "Pull latest changes from the branch you're currently" + [char]10 + "working on. "
$ws.Range("G5").Value = "Make sure to regulary push and commit to git"
$ws.Range("H5").Value = "No updates as of 08/05/2021"
$ws.Range("B5").WrapText = $true
$ws.Range("F5").WrapText = $true
$ws.Range("G5").WrapText = $true

# --- Row 6: Azure App Service Outage ---
$ws.Range("A6").Value = "Azure App Service Outage."
$ws.Range("B6").Value = "Client will not have access to the " + [char]10 + "production environment. "
$ws.Range("C6").Value = "low"
$ws.Range("D6").Value = "High"
$ws.Range("E6").Value = "Sean Kaila"
$ws.Range("F6").Value = "Deploy project on a parallel app service until the problem is resolved."
$ws.Range("G6").Value = "Regularly check server status " + [char]10 + "and set up a notification services to let somone know when the service has stopped working."
$ws.Range("H6").Value = "No updates as of 08/05/2021"
$ws.Range("B6").WrapText = $true
$ws.Range("F6").WrapText = $true
$ws.Range("G6").WrapText = $true

# --- Row 7: MySQL Database Faliure ---
$ws.Range("A7").Value = "MySQL Database Faliure"
$ws.Range("B7").Value = "Loss of core data stops clients from " + [char]10 + "using app. "
$ws.Range("C7").Value = "Low"
$ws.Range("D7").Value = "Meduim"
$ws.Range("E7").Value = "Sean Kaila"
$ws.Range("F7").Value = "Restore database from " + [char]10 + "a backup databse, so that clients go use app once again with minimal loss. "
$ws.Range("G7").Value = "Continue to back up Database twice a day" + [char]10 + "(Every 12 hours.) "
$ws.Range("H7").Value = "No updates as of 08/05/2021"
$ws.Range("B7").WrapText = $true
$ws.Range("F7").WrapText = $true
$ws.Range("G7").WrapText = $true

# --- Row 8: Edit feature failing ---
$ws.Range("A8").Value = "Edit feature failing "
$ws.Range("B8").Value = "Users not being able to use " + [char]10 + "cure functionality. "
$ws.Range("C8").Value = "Medium"
$ws.Range("D8").Value = "High"
$ws.Range("E8").Value = "Sean Kaila"
$ws.Range("F8").Value = "Investigate the issue as soon as " + [char]10 + "it arrives and deploy a hot fix."
$ws.Range("G8").Value = "Continue to Test edit functionality" + [char]10 + "making sure that its being tested correctly. "
$ws.Range("H8").Value = "No updates as of 08/05/2021"
$ws.Range("B8").WrapText = $true
$ws.Range("F8").WrapText = $true
$ws.Range("G8").WrapText = $true

# --- Row 9: Scope Creep ---
$ws.Range("A9").Value = "Scope Creep"
$ws.Range("B9").Value = "Results in core areas being unpolished " + [char]10 + "due to the amunt of taks needed to be complete in a short space of time."
$ws.Range("C9").Value = "high"
$ws.Range("D9").Value = "meduim"
$ws.Range("E9").Value = "Sean Kaila"
$ws.Range("F9").Value = "Delay the release of the product, " + [char]10 + "allowing for core functionality to be more polished. "
$ws.Range("G9").Value = "Revaluate the requirements against " + [char]10 + "what work is being done to minimuse the riskof spending too much time on unimportant tasks compered to the important ones."
$ws.Range("H9").Value = "No updates as of 08/05/2021"
$ws.Range("B9").WrapText = $true
$ws.Range("F9").WrapText = $true
$ws.Range("G9").WrapText = $true

Write-Host "values written"
